$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets hold the same event listing and need
# their attendance / price numbers refreshed to the latest scraped values.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: 南宁·三月三漫次元国风动漫节 - 想去人数 504 -> 509
    $ws.Range("F2").Value = 509

    # Row 3: 南宁·2024三月三国潮动漫节（良牙春典） - 想去人数 3395 -> 3407
    $ws.Range("F3").Value = 3407

    # Row 5: 南宁·布谷鸟动漫展4th - 想去人数 671 -> 674, 最低票价 35 -> 50
    $ws.Range("F5").Value = 674
    $ws.Range("G5").Value = 50
}
